$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 46056
$ws.Range("D8").Value = 158.47
$ws.Range("E8").Value = 149.37
$ws.Range("F8").Value = 159.37
$ws.Range("G8").Value = 149.26
$ws.Range("A9").Value = 46056
$ws.Range("D9").Value = 158.47
$ws.Range("E9").Value = 149.37
$ws.Range("F9").Value = 159.37
$ws.Range("G9").Value = 149.26
$ws.Range("A10").Value = 46056
$ws.Range("D10").Value = 159.88999999999999
$ws.Range("E10").Value = 150.88999999999999
$ws.Range("F10").Value = 160.88999999999999
$ws.Range("G10").Value = 151.13999999999999
$ws.Range("A11").Value = 46055
$ws.Range("D11").Value = 158.78
$ws.Range("E11").Value = 150.16999999999999
$ws.Range("F11").Value = 160.16999999999999
$ws.Range("G11").Value = 150.12
$ws.Range("A12").Value = 46055
$ws.Range("D12").Value = 158.78
$ws.Range("E12").Value = 150.16999999999999
$ws.Range("F12").Value = 160.16999999999999
$ws.Range("G12").Value = 150.12
$ws.Range("A13").Value = 46055
$ws.Range("D13").Value = 159.51
$ws.Range("E13").Value = 151.77000000000001
$ws.Range("F13").Value = 161.77000000000001
$ws.Range("G13").Value = 152.1
$ws.Range("A17").Value = 46056
$ws.Range("D17").Value = 163.47999999999999
$ws.Range("E17").Value = 154.35
$ws.Range("F17").Value = 164.35
$ws.Range("A18").Value = 46055
$ws.Range("D18").Value = 163.08000000000001
$ws.Range("E18").Value = 154.43
$ws.Range("F18").Value = 164.43
$ws.Range("A22").Value = 46056
$ws.Range("D22").Value = 159.54
$ws.Range("E22").Value = 150.83000000000001
$ws.Range("F22").Value = 160.43
$ws.Range("G22").Value = 152.58000000000001
$ws.Range("A23").Value = 46056
$ws.Range("D23").Value = 164.66
$ws.Range("E23").Value = 157.01
$ws.Range("F23").Value = 167.01
$ws.Range("A24").Value = 46056
$ws.Range("D24").Value = 164.82
$ws.Range("E24").Value = 157.65
$ws.Range("F24").Value = 167.65
$ws.Range("A25").Value = 46056
$ws.Range("D25").Value = 164.82
$ws.Range("E25").Value = 157.18
$ws.Range("F25").Value = 167.18
$ws.Range("G25").Value = 158.03
$ws.Range("A26").Value = 46056
$ws.Range("D26").Value = 164.4
$ws.Range("E26").Value = 158.76
$ws.Range("F26").Value = 168.76
$ws.Range("A27").Value = 46055
$ws.Range("D27").Value = 159.74
$ws.Range("E27").Value = 151.63
$ws.Range("F27").Value = 161.22999999999999
$ws.Range("G27").Value = 152.63
$ws.Range("A28").Value = 46055
$ws.Range("D28").Value = 164.28
$ws.Range("E28").Value = 157.12
$ws.Range("F28").Value = 167.12
$ws.Range("A29").Value = 46055
$ws.Range("D29").Value = 164.44
$ws.Range("E29").Value = 157.77000000000001
$ws.Range("F29").Value = 167.77
$ws.Range("A30").Value = 46055
$ws.Range("D30").Value = 164.44
$ws.Range("E30").Value = 157.30000000000001
$ws.Range("F30").Value = 167.3
$ws.Range("G30").Value = 157.35
$ws.Range("A31").Value = 46055
$ws.Range("D31").Value = 164.02
$ws.Range("E31").Value = 158.88999999999999
$ws.Range("F31").Value = 168.89
$ws.Range("A35").Value = 46056
$ws.Range("D35").Value = 158.82
$ws.Range("E35").Value = 148.69
$ws.Range("F35").Value = 157.69
$ws.Range("A36").Value = 46055
$ws.Range("D36").Value = 158.44
$ws.Range("E36").Value = 148.80000000000001
$ws.Range("F36").Value = 157.80000000000001
$ws.Range("A40").Value = 46056
$ws.Range("D40").Value = 164.2
$ws.Range("E40").Value = 156.07
$ws.Range("F40").Value = 166.07
$ws.Range("A41").Value = 46056
$ws.Range("D41").Value = 163.92
$ws.Range("E41").Value = 156.49
$ws.Range("F41").Value = 166.49
$ws.Range("A42").Value = 46055
$ws.Range("D42").Value = 163.84
$ws.Range("E42").Value = 156.05000000000001
$ws.Range("F42").Value = 166.05
$ws.Range("A43").Value = 46055
$ws.Range("D43").Value = 163.56
$ws.Range("E43").Value = 156.47
$ws.Range("F43").Value = 166.47
$ws.Range("A47").Value = 46056
$ws.Range("D47").Value = 157.99
$ws.Range("E47").Value = 150.68
$ws.Range("F47").Value = 160.68
$ws.Range("A48").Value = 46056
$ws.Range("D48").Value = 157.61000000000001
$ws.Range("E48").Value = 150.62
$ws.Range("F48").Value = 160.62
$ws.Range("A49").Value = 46055
$ws.Range("D49").Value = 158.47999999999999
$ws.Range("E49").Value = 150.94999999999999
$ws.Range("F49").Value = 160.94999999999999
$ws.Range("A50").Value = 46055
$ws.Range("D50").Value = 158.09
$ws.Range("E50").Value = 150.88
$ws.Range("F50").Value = 160.88
$ws.Range("A54").Value = 46056
$ws.Range("D54").Value = 173.19
$ws.Range("E54").Value = 164.46
$ws.Range("F54").Value = 174.46
$ws.Range("A55").Value = 46056
$ws.Range("D55").Value = 162.62
$ws.Range("E55").Value = 162.26
$ws.Range("F55").Value = 172.26
$ws.Range("A56").Value = 46056
$ws.Range("D56").Value = 162.5
$ws.Range("A57").Value = 46056
$ws.Range("D57").Value = 163.08000000000001
$ws.Range("E57").Value = 156.68
$ws.Range("A58").Value = 46056
$ws.Range("D58").Value = 158.85
$ws.Range("E58").Value = 152.58000000000001
$ws.Range("F58").Value = 162.58000000000001
$ws.Range("A59").Value = 46056
$ws.Range("D59").Value = 165.86
$ws.Range("E59").Value = 162.63999999999999
$ws.Range("A60").Value = 46055
$ws.Range("D60").Value = 172.83
$ws.Range("E60").Value = 164.65
$ws.Range("F60").Value = 174.65
$ws.Range("A61").Value = 46055
$ws.Range("D61").Value = 165.44
$ws.Range("E61").Value = 163.12
$ws.Range("F61").Value = 173.12
$ws.Range("A62").Value = 46055
$ws.Range("D62").Value = 162.46
$ws.Range("A63").Value = 46055
$ws.Range("D63").Value = 163
$ws.Range("E63").Value = 157.54
$ws.Range("A64").Value = 46055
$ws.Range("D64").Value = 158.77000000000001
$ws.Range("E64").Value = 153.44
$ws.Range("F64").Value = 163.44
$ws.Range("A65").Value = 46055
$ws.Range("D65").Value = 165.45
$ws.Range("E65").Value = 162.80000000000001

$ws.PageSetup.Zoom = 67

